$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$ws1.Range("G3").Value = "2016-08-28 22:48:10"

$ws2.Range("H3").Value = "2016-08-28 22:48:06"
$ws2.Range("K3").Value = "2016-08-28 22:48:22"

$ws3.Range("H3").Value = "2016-08-28 22:48:10"
$ws3.Range("K3").Value = "2016-08-28 22:48:29"
